$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-9: Columns H (Correct Answer) and I (Time in seconds) were stored
# as text (inlineStr) but should really be numeric values.
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 45

$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 60

$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 45

$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 60

$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 45

$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 30

$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 45

$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 30

# Row 10: answer text updated, correct-answer column cleared, time made numeric.
# Leading apostrophe forces these to stay literal text (otherwise Excel would
# auto-convert "1.0" to the number 1, and "" simply clears the cell).
$ws.Range("C10").Value = "'1.0"
$ws.Range("H10").Value = "'"
$ws.Range("I10").Value = 45

# Row 11: time made numeric
$ws.Range("I11").Value = 60
